$wb = $excel.ActiveWorkbook

# --- APT_ATFM_SES_YY (sheet1): update C/D raw values; B (formula) recalculates ---
$ws1 = $wb.Worksheets.Item("APT_ATFM_SES_YY")
$ws1.Range("C6").Value = 2455264
$ws1.Range("D6").Value = 1348999
$ws1.Range("C7").Value = 2540608
$ws1.Range("D7").Value = 1484416
$ws1.Range("C8").Value = 2608155
$ws1.Range("D8").Value = 1753480
$ws1.Range("C9").Value = 2675793
$ws1.Range("D9").Value = 2234754
$ws1.Range("C10").Value = 1216967
$ws1.Range("D10").Value = 578133

# --- APT_ATFM_SES_MM (sheet2): update C/D raw values; B/E (formulas) recalculate ---
$ws2 = $wb.Worksheets.Item("APT_ATFM_SES_MM")
$ws2.Range("C6").Value = 397807
$ws2.Range("D6").Value = 188850
$ws2.Range("C7").Value = 375556
$ws2.Range("D7").Value = 168593
$ws2.Range("C8").Value = 430919
$ws2.Range("D8").Value = 437216
$ws2.Range("C9").Value = 463047
$ws2.Range("D9").Value = 336793
$ws2.Range("C10").Value = 497647
$ws2.Range("D10").Value = 520464
$ws2.Range("C11").Value = 510817
$ws2.Range("D11").Value = 582838
$ws2.Range("C12").Value = 531594
$ws2.Range("D12").Value = 563182
$ws2.Range("C13").Value = 519378
$ws2.Range("C14").Value = 512430
$ws2.Range("D14").Value = 490127
$ws2.Range("C15").Value = 492916
$ws2.Range("D15").Value = 385218
$ws2.Range("C16").Value = 406334
$ws2.Range("D16").Value = 183427
$ws2.Range("C17").Value = 401805
$ws2.Range("D17").Value = 291252
$ws2.Range("C18").Value = 395118
$ws2.Range("D18").Value = 237888
$ws2.Range("C19").Value = 378978
$ws2.Range("C20").Value = 245513
$ws2.Range("C21").Value = 45174
$ws2.Range("C22").Value = 57679
$ws2.Range("C23").Value = 94505
$ws2.Range("C24").Value = 205959
$ws2.Range("D24").Value = 13893
$ws2.Range("C25").Value = 249406
$ws2.Range("D25").Value = 16688
$ws2.Range("C26").Value = 220299
$ws2.Range("C27").Value = 194469
$ws2.Range("D27").Value = 10044
$ws2.Range("C28").Value = 132028

# --- Update cell selections (per-sheet cursor position) ---
# APT_ATFM_SES_YY: selection moves to D39
$ws1.Range("D39").Select()

# APT_ATFM_LOC (sheet3): selection moves from H18 to E8.
# Select it LAST so this sheet remains the workbook's active tab,
# matching the original file (activeTab = APT_ATFM_LOC / tabSelected).
$ws3 = $wb.Worksheets.Item("APT_ATFM_LOC")
$ws3.Range("E8").Select()

Write-Host "edits applied"
